$d = $word.ActiveDocument

# --- 1) Strike through "Progress bar" and "Auto-decrement/increment" bullets ---
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd("`r", "`n")
    if ($t -eq "Progress bar" -or $t -eq "Auto-decrement/increment") {
        $p.Range.Font.StrikeThrough = $true
    }
}

# --- 2) Append three new list paragraphs after the final paragraph ---
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$newp = $d.Paragraphs.Item($d.Paragraphs.Count)
$rInsert = $newp.Range
$rInsert.Collapse(0)

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>General Implementation</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Java 23</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>GUI using Swing or JavaFX</w:t></w:r></w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rInsert.InsertXML($xmlFrag)

# InsertXML leaves a stray empty paragraph (carrying the old trailing pPr) after
# the inserted content; merge it away so the document ends cleanly.
$trailing = $d.Paragraphs.Item($d.Paragraphs.Count)
$delRange = $d.Range($trailing.Range.Start - 1, $trailing.Range.End)
$delRange.Delete()

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
